# CIV-11239 fix full stop
#
# The "Costs reserved" option text in the order template has a stray
# full stop immediately before the <<costsReservedText>> merge field,
# producing a double full-stop once that field is populated
# (e.g. "Costs reserved. Some reason text."). Remove the stray period,
# keeping the single space that separates "Costs reserved" from the
# merge field.

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "Costs reserved. ",   # FindText
    $true,                 # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                 # Forward
    1,                      # Wrap (wdFindContinue)
    $false,                # Format
    "Costs reserved ",    # ReplaceWith
    2                       # Replace (wdReplaceAll)
)

Write-Output "Costs reserved fix applied: $found"
